$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.481.65"
$ws.Range("E2").Value = "  -0.23%  "

$ws.Range("D3").Value = "2.293.62"
$ws.Range("E3").Value = "  -1.59%  "

$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").Value = "'540.61"
$ws.Range("E5").Value = "  -0.81%  "

$ws.Range("D6").Value = "'128.64"
$ws.Range("E6").Value = "  -2.95%  "

$ws.Range("E7").Value = "  +0.15%  "

$ws.Range("D8").Value = "'0.569"
$ws.Range("E8").Value = "  -2.79%  "

$ws.Range("D9").Value = "2.292.28"
$ws.Range("E9").Value = "  -1.50%  "

$ws.Range("D10").Value = "'0.100"
$ws.Range("E10").Value = "  -0.19%  "

$ws.Range("E11").Value = "  -0.33%  "

$ws.Range("D12").Value = "'0.150"
$ws.Range("E12").Value = "  -0.55%  "

$ws.Range("D13").Value = "'0.330"
$ws.Range("E13").Value = "  -1.50%  "

$ws.Range("D14").Value = "'23.07"
$ws.Range("E14").Value = "  -3.45%  "

$ws.Range("D15").Value = "2.703.11"
$ws.Range("E15").Value = "  -1.23%  "

$ws.Range("D16").Value = "59.421.70"
$ws.Range("E16").Value = "  +0.12%  "

$ws.Range("E17").Value = "  -1.64%  "

$ws.Range("D18").Value = "2.293.63"
$ws.Range("E18").Value = "  -0.76%  "

$ws.Range("D19").Value = "'10.39"
$ws.Range("E19").Value = "  -2.25%  "

$ws.Range("E20").Value = "  -4.15%  "

$ws.Range("D21").Value = "'309.09"
$ws.Range("E21").Value = "  -1.89%  "

$ws.Range("D22").Value = "'6.49"
$ws.Range("E22").Value = "  -2.67%  "

$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  -0.21%  "

$ws.Range("D24").Value = "'62.70"
$ws.Range("E24").Value = "  -0.29%  "

$ws.Range("D25").Value = "'0.167"
$ws.Range("E25").Value = "  -3.83%  "

$ws.Range("E26").Value = "  +0.43%  "

$ws.Range("D27").Value = "'7.69"
$ws.Range("E27").Value = "  -3.82%  "

$ws.Range("D28").Value = "'1.34"
$ws.Range("E28").Value = "  +2.04%  "

$ws.Range("E29").Value = "  +1.30%  "

$ws.Range("D30").Value = "'171.80"
$ws.Range("E30").Value = "  +0.39%  "

$ws.Range("E31").Value = "  -1.56%  "

$ws.Range("D32").Value = "0.0₃0716"
$ws.Range("E32").Value = "  -3.63%  "

$ws.Range("D33").Value = "'5.77"
$ws.Range("E33").Value = "  -2.36%  "

$ws.Range("D34").Value = "'0.376"
$ws.Range("E34").Value = "  -2.28%  "

$ws.Range("E36").Value = "  -6.29%  "

$ws.Range("E37").Value = "  -1.83%  "

$ws.Range("E38").Value = "  +0.09%  "

$ws.Range("D39").Value = "'3.97"
$ws.Range("E39").Value = "  -2.88%  "

$ws.Range("D40").Value = "'311.05"
$ws.Range("E40").Value = "  -2.36%  "

$ws.Range("D41").Value = "'37.38"
$ws.Range("E41").Value = "  -1.85%  "

$ws.Range("E42").Value = "  -2.65%  "

$ws.Range("D43").Value = "'136.20"
$ws.Range("E43").Value = "  -5.27%  "

$ws.Range("D44").Value = "'3.41"
$ws.Range("E44").Value = "  -1.08%  "

$ws.Range("E45").Value = "  -1.65%  "

$ws.Range("E46").Value = "  -0.14%  "

$ws.Range("D47").Value = "'18.43"
$ws.Range("E47").Value = "  +0.01%  "

$ws.Range("D48").Value = "'0.0486"
$ws.Range("E48").Value = "  -2.14%  "

$ws.Range("D49").Value = "0.0₆0223"
$ws.Range("E49").Value = "  +22.74%  "

$ws.Range("D50").Value = "'0.0212"
$ws.Range("E50").Value = "  +0.57%  "

$ws.Range("D51").Value = "'10.98"
$ws.Range("E51").Value = "  -0.33%  "
